# aggiornamento fino a 1/09/2021
# Append new daily rows (358-366) to the "Bomporto" report sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: row, date-serial (A), nuovi pos. (B), somma mobile 7gg. (C), somma mobile 7gg. per 100mila abitanti (D)
$rows = @(
    @(358, 44432, 0, 7,  69.51340615690168),
    @(359, 44433, 1, 8,  79.44389275074478),
    @(360, 44434, 0, 7,  69.51340615690168),
    @(361, 44435, 2, 5,  49.65243296921549),
    @(362, 44436, 0, 4,  39.72194637537239),
    @(363, 44437, 7, 11, 109.2353525322741),
    @(364, 44438, 2, 12, 119.1658391261172),
    @(365, 44439, 1, 13, 129.0963257199603),
    @(366, 44440, 1, 13, 129.0963257199603)
)

# Use the last existing data row (357) as the formatting template for the new rows.
$templateRow = 357

foreach ($r in $rows) {
    $rowIndex = $r[0]

    $ws.Cells.Item($rowIndex, 1).Value2 = $r[1]
    $ws.Cells.Item($rowIndex, 2).Value2 = $r[2]
    $ws.Cells.Item($rowIndex, 3).Value2 = $r[3]
    $ws.Cells.Item($rowIndex, 4).Value2 = $r[4]

    # Copy formatting (number format, style, borders, alignment) from the template row
    $srcRow = $ws.Range("A" + $templateRow + ":D" + $templateRow)
    $dstRow = $ws.Range("A" + $rowIndex + ":D" + $rowIndex)
    $srcRow.Copy()
    $dstRow.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
